# Applies the "Added my part to the progress report" edit:
# rewrites Team Member #1's three answers and relocates the _GoBack bookmark.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Team Member #1 - "What I did since last progress report:"
# ---------------------------------------------------------------------
$old1 = "Since the last project report, I have worked on getting database tables and diagrams set up for use on the prototype. Relationships between the tables have also been established. I contacted John Nordlie about having a SQL Server database hosted on the UND campus, and we are currently corresponding about this."
$new1 = "Since the last project report, I have contacted John Nordlie and discussed with him our options for hosting a database on a UND server. Because of the move from Streibel to the engineering buildings, there will not be a chance to get anything running on campus. This will result in us taking a major shift into the way we will be developing the inventory system. SQL Server may no longer be an option for us, and it had been determined that we will likely have to change our approach."
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------
# Team Member #1 - "Anything holding me up..."
# ---------------------------------------------------------------------
$old2 = "Homework for other classes has still been keeping me pretty busy, but my weekends should be more free in these next few weeks. Correspondence with John Nordlie will also be limiting us slightly while we wait for instruction on how to continue moving forward with our database."
$new2 = "Final projects and homework assignments have been keeping me very busy in the past few weeks. After the final exam week I will be able to make some progress in implementing our database and work towards a final decision on which software we will be using. Another setback we had was learning that we won" + [char]0x2019 + "t be able to host a database on campus."
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------
# Team Member #1 - "What I plan to do before the next progress report:"
# ---------------------------------------------------------------------
$old3 = "Before our next progress report, I would like to have our database hosted on the UND servers, but this is not a major priority. It is possible that we won't be able to have something up and running before the next report, but we can adjust our work as needed. I also want to try to get the current spreadsheet from the Computer Science department so I can make sure our database tables and fields contain the necessary information, and so I can import the data for testing. Lastly, I want to confirm which fields are necessary, and if we are lacking data for those fields, how to go about fixing the data to fit our needs."
$new3 = "Before our next progress report, I would like to discuss with my team members which direction we want to go with our project after these new developments. We are debating between continuing with what we have already, starting new with a .NET website, or starting new using mainly JavaScript. Database software may also change for us to either MySQL or SQLite. We will make a decision before the next progress report."
$d.Content.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# ---------------------------------------------------------------------
# Relocate the _GoBack bookmark from the end of the document to right
# after Team Member #1's first answer (where the last edit happened).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$anchor = $d.Content.Find.Execute("it had been determined that we will likely have to change our approach.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Content
$target.Start = $target.End
$target.Collapse(0)
$found = $d.Content
$found.Find.Execute("it had been determined that we will likely have to change our approach.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$found.Collapse(0)
$d.Bookmarks.Add("_GoBack", $found) | Out-Null
